# "bonus gera + bonus M1"
# - rename header "executivos" -> "executivo"
# - bump the "m1" (col E) bonus numbers for Marcela, Eder and Tatiana,
#   which ripples into the "ponto gera" (col I) formula results
# - widen columns C:I to a uniform width (was split 4/9/9.66)
# - leave the selection on E5 (where the last edit, E6's neighbour, was made)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "executivo"

$ws.Range("E2").Value = 9
$ws.Range("E3").Value = 9
$ws.Range("E6").Value = 5

$ws.Range("C1:I1").ColumnWidth = 9.94

$ws.Range("E5").Select()
